$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.061.64"
$ws.Range("E2").Value = "  -0.22%  "

$ws.Range("D3").Value = "3.862.83"
$ws.Range("E3").Value = "  +1.38%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "701.24"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.65"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.23%  "

$ws.Range("E7").Value = "  +1.34%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("E10").Value = "  -0.27%  "

$ws.Range("E11").Value = "  -4.30%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000259"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.34"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.12%  "

$ws.Range("D15").Value = "4.513.45"
$ws.Range("E15").Value = "  +1.48%  "

$ws.Range("D16").Value = "3.982.37"
$ws.Range("E16").Value = "  +4.47%  "

$ws.Range("D17").Value = "71.087.50"
$ws.Range("E17").Value = "  -0.16%  "

$ws.Range("E18").Value = "  -0.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.46"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.32%  "

$ws.Range("E20").Value = "  -0.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "500.78"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.74"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.34%  "

$ws.Range("E23").Value = "  +0.77%  "

$ws.Range("E24").Value = "  +2.33%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.07"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.69"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.23"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.05%  "

$ws.Range("E28").Value = "  -2.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.17"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("E31").Value = "  -0.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.28"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.57%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.64"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.182"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.82%  "

$ws.Range("E35").Value = "  -0.15%  "

$ws.Range("D36").Value = "3.819.11"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.07%  "

$ws.Range("E38").Value = "  +1.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.41"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +8.71%  "

$ws.Range("E40").Value = "  +8.93%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.44"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.73%  "

$ws.Range("E42").Value = "  +0.77%  "

$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("E44").Value = "  +0.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000317"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.35%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.86"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.19%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "49.20"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "418.17"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.87%  "

$ws.Range("E49").Value = "  +0.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.58"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.36%  "

